$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Thu Nov 07 16:46:45 EST 2024"
$ws.Range("B3").Value = "Thu Nov 07 16:47:00 EST 2024"
$ws.Range("B4").Value = "Thu Nov 07 16:47:14 EST 2024"
